$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2 through 27 holds a date serial
# that was bumped forward by one day: 45188 (2023-09-19) -> 45189 (2023-09-20).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45189
}
